# Feat: added input validation to loadEnsembleStructure and respective tests
#
# The underlying change to the "metsData" worksheet is the removal of two
# metabolites (m4 and m22) that no longer pass the new input validation.
# Removing them shifts every subsequent metabolite row up, while each
# metabolite keeps its own bound values (columns B/C/D).
#
# In the original (pre-edit) row numbering:
#   row 5  -> m4   (A5 = "m4")
#   row 23 -> m22  (A23 = "m22")
# Delete from the bottom up so row indices for the earlier deletion are
# not invalidated by the later one.

$wb = $excel.ActiveWorkbook

$metsData = $wb.Worksheets.Item("metsData")
$metsData.Rows.Item(23).EntireRow.Delete()
$metsData.Rows.Item(5).EntireRow.Delete()

# Make metsData the active/selected sheet with A5 selected (matches the
# sheetView/selection recorded for the sheet after the edit), which also
# flips the previously active "general" sheet's tabSelected flag off and
# updates the workbook's activeTab index automatically.
$metsData.Activate()
$metsData.Range("A5").Select()

# Trailing near-empty rows that show up at the very bottom of the sheet
# after the row deletions.
$metsData.Rows.Item(1048575).RowHeight = 12.8
$metsData.Rows.Item(1048576).RowHeight = 12.8

# Minor column width adjustments (auto-fit sized columns) on a few other
# sheets that were recalculated as part of the same save.
$general = $wb.Worksheets.Item("general")
$general.Columns.Item(1).ColumnWidth = 58.333333333333336

$kinetics1 = $wb.Worksheets.Item("kinetics1")
$kinetics1.Columns.Item(2).ColumnWidth = 12.333333333333334
$kinetics1.Columns.Item(3).ColumnWidth = 34.666666666666664
$kinetics1.Columns.Item(4).ColumnWidth = 35.666666666666664

$mets = $wb.Worksheets.Item("mets")
$mets.Columns.Item(1).ColumnWidth = 27.5
$mets.Columns.Item(2).ColumnWidth = 15.5
